$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update Marking row (row 11) - Right marks value
$ws.Range("B11").Value = 5

# Update Total row (row 12) - Right marks total
$ws.Range("B12").Value = 50

# Update Max/score summary text
$ws.Range("E12").Value = "50/140"
